$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks; we will re-add them (shifted by one row) after
# inserting the new row, since Insert() does not re-target hyperlink refs.
$ws.Hyperlinks.Delete()

# Insert a brand new row 2 (pushes the current rows 2-13 down to 3-14).
$ws.Rows("2:2").Insert()

# The freshly inserted row inherits the header row's formatting; clear it so
# the new row looks like an ordinary data row, and drop the placeholder
# empty cells in columns B, C, E and F (they must not be present at all).
$ws.Range("B2:C2").Clear()
$ws.Range("E2:F2").Clear()
$ws.Range("A2").ClearFormats()
$ws.Range("D2").ClearFormats()
# G2 should use the same "URL" style as the other link cells (e.g. G3).
$ws.Range("G2").Style = $ws.Range("G3").Style

# Populate the new row. Write the URL (G) before the name (A) so that the
# shared-string table gains the two new entries in the same order as the
# source edit (URL first, then the game name).
$ws.Range("G2").Value = "https://pan.baidu.com/s/1wlbvZ3w_KPL3bqt7nNtQCg?pwd=rftt"
$ws.Range("D2").Value = "rftt"
$ws.Range("A2").Value = "宝可梦朱紫"

# Re-create the hyperlinks on column G, each shifted down one row from where
# it used to be (row N -> row N+1), preserving the special "list/path=%2F"
# sub-address + display text on the two entries that had them. Adding a
# TextToDisplay replaces the anchor cell's text, so immediately restore the
# cell's real (shorter) URL text afterwards for those two special cases.
$h1 = $ws.Hyperlinks.Add($ws.Range("G4"), "https://pan.baidu.com/s/1x_V0cQZyzhAIzr97GCznlA?pwd=5fub", "list/path=%2F")
$h1.TextToDisplay = "https://pan.baidu.com/s/1x_V0cQZyzhAIzr97GCznlA?pwd=5fub - list/path=%2F"
$ws.Range("G4").Value = "https://pan.baidu.com/s/1x_V0cQZyzhAIzr97GCznlA?pwd=5fub#list/path=%2F"

$ws.Hyperlinks.Add($ws.Range("G12"), "https://pan.baidu.com/s/1ZbWggC3GDJv7BUgxTIbGzg")
$ws.Hyperlinks.Add($ws.Range("G3"), "https://pan.baidu.com/s/1KykYnfqctZOEDgJp_nxGsA?pwd=uqer")
$ws.Hyperlinks.Add($ws.Range("G11"), "https://pan.baidu.com/s/1zKgW1pjqUnZ2dtEq2xXFMw?pwd=ccx6")
$ws.Hyperlinks.Add($ws.Range("G5"), "https://pan.baidu.com/s/1UARljz8BQP1uTU3Lie_2oQ")

$h2 = $ws.Hyperlinks.Add($ws.Range("G14"), "https://pan.baidu.com/s/1Re4OiBosRO_y77sDJRBRuw", "list/path=%2F")
$h2.TextToDisplay = "https://pan.baidu.com/s/1Re4OiBosRO_y77sDJRBRuw - list/path=%2F"
$ws.Range("G14").Value = "https://pan.baidu.com/s/1Re4OiBosRO_y77sDJRBRuw#list/path=%2F"

$ws.Hyperlinks.Add($ws.Range("G9"), "https://pan.baidu.com/share/init?surl=n6ivaYdevwiyNpXc1Fgpxg")
$ws.Hyperlinks.Add($ws.Range("G7"), "https://pan.baidu.com/s/1zGN-60z0e43PEJYohxM8cA")
$ws.Hyperlinks.Add($ws.Range("G10"), "https://pan.baidu.com/s/1sF9Suvp0mVQNP56NZmL64A")
$ws.Hyperlinks.Add($ws.Range("G8"), "https://pan.baidu.com/s/1EWpU8lPT_bYMm3uoK174kA?pwd=5vch")
$ws.Hyperlinks.Add($ws.Range("G6"), "https://pan.baidu.com/s/1cIIZZv89eBKv255fKtD4cQ")

# Match the author's final cursor position.
$ws.Range("F18").Select()
